$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43; existing rows 43-72 shift down to 44-73.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new price observation (same categorical
# context as its neighbours, new date / volume / price data).
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44729
$ws.Range("D43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100102
$ws.Range("H43").Value = "Cítricos"
$ws.Range("I43").Value = 100102006
$ws.Range("J43").Value = "Pomelo"
$ws.Range("K43").Value = "Start Ruby"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 300
$ws.Range("N43").Value = 8000
$ws.Range("O43").Value = 8000
$ws.Range("P43").Value = 8000
$ws.Range("Q43").Value = "$/caja 14 kilos"
$ws.Range("R43").Value = "Región Metropolitana"
$ws.Range("S43").Value = 571
$ws.Range("T43").Value = 14
